# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stocks) sheet gets three new trailing columns added:
#   H = date             (e.g. "2012-04-19")
#   I = legislator_name  (e.g. "李慶華")
#   J = legislator_id    (e.g. 607)
# for every existing data row, mirroring the same three columns already
# present on the other sheets in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Pull formats from the last existing header/data cell (column G) so the
# new columns look identical to the rest of the table (bold+border header,
# plain data row) instead of defaulting to unformatted cells.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G2").Copy()
$ws.Range("H2:J2").PasteSpecial(-4122)   # xlPasteFormats

# New header row.
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# New data row. Force the date column to text first so "2012-04-19" is
# stored verbatim instead of being auto-converted into a date serial
# number.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2012-04-19"
$ws.Range("I2").Value = "李慶華"
$ws.Range("J2").Value = 607
